$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.771.86"
$ws.Range("E2").Value = "  -4.29%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.653.91"
$ws.Range("E3").Value = "  -6.56%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.007"
$ws.Range("E4").Value = "  +0.73%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.07"
$ws.Range("E5").Value = "  -2.76%  "
$ws.Range("E6").Value = "  +0.67%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3643"
$ws.Range("E7").Value = "  -5.22%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3313"
$ws.Range("E8").Value = "  -9.46%  "
$ws.Range("E9").Value = "  -8.09%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.144"
$ws.Range("E10").Value = "  -7.37%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07190"
$ws.Range("E11").Value = "  -6.44%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.004"
$ws.Range("E12").Value = "  +0.81%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.094"
$ws.Range("E13").Value = "  -6.08%  "
$ws.Range("E14").Value = "  -8.40%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.711"
$ws.Range("E15").Value = "  -6.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.651.36"
$ws.Range("E16").Value = "  -6.48%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001073"
$ws.Range("E17").Value = "  -7.74%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06603"
$ws.Range("E18").Value = "  -3.79%  "
$ws.Range("E19").Value = "  +0.57%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "80.54"
$ws.Range("E20").Value = "  -7.91%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.50"
$ws.Range("E21").Value = "  -7.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.056"
$ws.Range("E22").Value = "  -7.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.26"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.761.28"
$ws.Range("E24").Value = "  -4.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.418"
$ws.Range("E25").Value = "  -0.37%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.576"
$ws.Range("E26").Value = "  -13.28%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "150.18"
$ws.Range("E27").Value = "  -3.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.42"
$ws.Range("E28").Value = "  -6.69%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "128.29"
$ws.Range("E29").Value = "  -4.79%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.839.56"
$ws.Range("E30").Value = "  -6.26%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.209"
$ws.Range("E31").Value = "  +0.22%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.133"
$ws.Range("E32").Value = "  -3.80%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.216"
$ws.Range("E33").Value = "  -13.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.722"
$ws.Range("E34").Value = "  -5.27%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08502"
$ws.Range("E35").Value = "  -2.56%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.96"
$ws.Range("E36").Value = "  -8.45%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.286"
$ws.Range("E37").Value = "  -7.42%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06279"
$ws.Range("E38").Value = "  -7.47%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02304"
$ws.Range("E39").Value = "  -7.57%  "
$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2118"
$ws.Range("E40").Value = "  -5.26%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.431"
$ws.Range("E41").Value = "  -10.24%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.229"
$ws.Range("E42").Value = "  -5.94%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6142"
$ws.Range("E43").Value = "  -6.97%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  +0.56%  "
$ws.Range("E45").Value = "  -6.64%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.766"
$ws.Range("E46").Value = "  -4.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5842"
$ws.Range("E47").Value = "  -8.47%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.006"
$ws.Range("E48").Value = "  -7.99%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "124.08"
$ws.Range("E49").Value = "  -6.68%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07078"
$ws.Range("E50").Value = "  -5.74%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "75.52"
$ws.Range("E51").Value = "  -6.75%  "
